$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2, shifting existing data rows (2-10) down to (3-11)
$ws.Rows.Item(2).Insert()

# Populate the new row 2 with the new weekly record
$ws.Cells.Item(2, 1).Value = 11
$ws.Cells.Item(2, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(2, 3).Value = "Bíobío"
$ws.Cells.Item(2, 4).Value = 44545
$ws.Cells.Item(2, 4).Style = $ws.Cells.Item(3, 4).Style
$ws.Cells.Item(2, 4).NumberFormat = $ws.Cells.Item(3, 4).NumberFormat
$ws.Cells.Item(2, 5).Value = 8
$ws.Cells.Item(2, 6).Value = 300000000
$ws.Cells.Item(2, 7).Value = "Espárragos"
$ws.Cells.Item(2, 8).Value = "Sin especificar"
$ws.Cells.Item(2, 9).Value = "Primera"
$ws.Cells.Item(2, 10).Value = 550
$ws.Cells.Item(2, 11).Value = 1700
$ws.Cells.Item(2, 12).Value = 1800
$ws.Cells.Item(2, 13).Value = 1755
$ws.Cells.Item(2, 14).Value = "`$/kilo"
$ws.Cells.Item(2, 15).Value = "Provincia de Linares"
$ws.Cells.Item(2, 16).Value = 1755
$ws.Cells.Item(2, 17).Value = 1
$ws.Cells.Item(2, 18).Value = "Hortaliza"
